$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 328.42856
$ws.Range("I12").Value = 354.83334
$ws.Range("K12").Value = 354.83334
$ws.Range("M12").Value = -184.83334
$ws.Range("H33").Value = 1427.6
$ws.Range("I33").Value = 1154.8636
$ws.Range("J33").Value = 2177.625
$ws.Range("K33").Value = 1154.8636
$ws.Range("L33").Value = 2177.625
$ws.Range("M33").Value = -925.8635999999999
$ws.Range("N33").Value = -2635.625
$ws.Range("H52").Value = 8665.666999999999
$ws.Range("H62").Value = 2971.25
$ws.Range("I62").Value = 2294
$ws.Range("K62").Value = 2294
$ws.Range("M62").Value = -1670
$ws.Range("H65").Value = 2971.25
$ws.Range("I65").Value = 2294
$ws.Range("K65").Value = 11470
$ws.Range("M65").Value = -8350
$ws.Range("H74").Value = 4807.615
$ws.Range("H76").Value = 4282.3335
$ws.Range("I76").Value = 3827
$ws.Range("J76").Value = 4919.8
$ws.Range("K76").Value = 3827
$ws.Range("L76").Value = 4919.8
$ws.Range("M76").Value = -3512
$ws.Range("N76").Value = -5549.8
$ws.Range("H77").Value = 4807.615
$ws.Range("H79").Value = 4282.3335
$ws.Range("I79").Value = 3827
$ws.Range("J79").Value = 4919.8
$ws.Range("K79").Value = 3827
$ws.Range("L79").Value = 4919.8
$ws.Range("M79").Value = -2735
$ws.Range("N79").Value = -7103.8
$ws.Range("H138").Value = 2116.7627
$ws.Range("J138").Value = 2593.4285
$ws.Range("L138").Value = 7780.2855
$ws.Range("N138").Value = -18060.2855

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 93332
$ws.Range("J24").Value = 93332
$ws.Range("L24").Value = 93332
$ws.Range("N24").Value = -94080
$ws.Range("H32").Value = 10418854
$ws.Range("I32").Value = 10418854
$ws.Range("K32").Value = 10418854
$ws.Range("M32").Value = -10418567
$ws.Range("H45").Value = 1856.6522
$ws.Range("I45").Value = 1294.4667
$ws.Range("J45").Value = 2910.75
$ws.Range("K45").Value = 1294.4667
$ws.Range("L45").Value = 2910.75
$ws.Range("M45").Value = -917.4666999999999
$ws.Range("N45").Value = -3664.75
$ws.Range("H61").Value = 8794510
$ws.Range("I61").Value = 10420647
$ws.Range("K61").Value = 10420647
$ws.Range("M61").Value = -10420435
$ws.Range("H88").Value = 2628.25
$ws.Range("I88").Value = 2753
$ws.Range("J88").Value = 2503.5
$ws.Range("K88").Value = 2753
$ws.Range("L88").Value = 2503.5
$ws.Range("M88").Value = -2347
$ws.Range("N88").Value = -3315.5
$ws.Range("H91").Value = 2628.25
$ws.Range("I91").Value = 2753
$ws.Range("J91").Value = 2503.5
$ws.Range("K91").Value = 2753
$ws.Range("L91").Value = 2503.5
$ws.Range("M91").Value = -1349
$ws.Range("N91").Value = -5311.5
$ws.Range("H100").Value = 93332
$ws.Range("J100").Value = 93332
$ws.Range("L100").Value = 93332
$ws.Range("N100").Value = -95496
$ws.Range("H102").Value = 10433.392
$ws.Range("I102").Value = 11303.238
$ws.Range("J102").Value = 1300
$ws.Range("K102").Value = 11303.238
$ws.Range("L102").Value = 1300
$ws.Range("M102").Value = -9681.237999999999
$ws.Range("N102").Value = -4544
$ws.Range("H103").Value = 49500
$ws.Range("J103").Value = 49500
$ws.Range("L103").Value = 49500
$ws.Range("N103").Value = -51844
$ws.Range("H120").Value = 69996.664
$ws.Range("J120").Value = 69996.664
$ws.Range("L120").Value = 69996.664
$ws.Range("N120").Value = -79672.664
$ws.Range("H136").Value = 8794510
$ws.Range("I136").Value = 10420647
$ws.Range("K136").Value = 31261941
$ws.Range("M136").Value = -31259391

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1401.3572
$ws.Range("I80").Value = 1176.75
$ws.Range("J80").Value = 1491.2
$ws.Range("K80").Value = 1176.75
$ws.Range("L80").Value = 1491.2
$ws.Range("M80").Value = -178.75
$ws.Range("N80").Value = -3487.2
$ws.Range("H83").Value = 1401.3572
$ws.Range("I83").Value = 1176.75
$ws.Range("J83").Value = 1491.2
$ws.Range("K83").Value = 5883.75
$ws.Range("L83").Value = 7456
$ws.Range("M83").Value = -891.75
$ws.Range("N83").Value = -17440
$ws.Range("H86").Value = 2733.9285
$ws.Range("I86").Value = 2462.3635
$ws.Range("K86").Value = 2462.3635
$ws.Range("M86").Value = -1339.3635
$ws.Range("H89").Value = 2733.9285
$ws.Range("I89").Value = 2462.3635
$ws.Range("K89").Value = 12311.8175
$ws.Range("M89").Value = -6695.817499999999
$ws.Range("H94").Value = 2004.7
$ws.Range("I94").Value = 2256.125
$ws.Range("J94").Value = 999
$ws.Range("K94").Value = 2256.125
$ws.Range("L94").Value = 999
$ws.Range("M94").Value = -1805.125
$ws.Range("N94").Value = -1901
$ws.Range("H98").Value = 77860.336
$ws.Range("J98").Value = 77860.336
$ws.Range("L98").Value = 77860.336
$ws.Range("N98").Value = -83850.336

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 586513.1
$ws.Range("I31").Value = 3831.4707
$ws.Range("J31").Value = 3888376
$ws.Range("K31").Value = 3831.4707
$ws.Range("L31").Value = 3888376
$ws.Range("M31").Value = -3536.4707
$ws.Range("N31").Value = -3888966
$ws.Range("H34").Value = 586513.1
$ws.Range("I34").Value = 3831.4707
$ws.Range("J34").Value = 3888376
$ws.Range("K34").Value = 3831.4707
$ws.Range("L34").Value = 3888376
$ws.Range("M34").Value = -3629.4707
$ws.Range("N34").Value = -3888780
$ws.Range("H105").Value = 2228.4666
$ws.Range("I105").Value = 1785.75
$ws.Range("K105").Value = 1785.75
$ws.Range("M105").Value = -38.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 4464.636
$ws.Range("I137").Value = 6994
$ws.Range("J137").Value = 3516.125
$ws.Range("K137").Value = 20982
$ws.Range("L137").Value = 10548.375
$ws.Range("M137").Value = -15882
$ws.Range("N137").Value = -20748.375

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 443.85715
$ws.Range("I2").Value = 68.625
$ws.Range("K2").Value = 68.625
$ws.Range("M2").Value = 44.375
$ws.Range("H63").Value = 15214
$ws.Range("J63").Value = 15214
$ws.Range("L63").Value = 15214
$ws.Range("N63").Value = -16586
$ws.Range("H66").Value = 15214
$ws.Range("J66").Value = 15214
$ws.Range("L66").Value = 45642
$ws.Range("N66").Value = -52506
$ws.Range("H80").Value = 3309.6
$ws.Range("I80").Value = 3437
$ws.Range("J80").Value = 2800
$ws.Range("K80").Value = 3437
$ws.Range("L80").Value = 2800
$ws.Range("M80").Value = -2439
$ws.Range("N80").Value = -4796
$ws.Range("H83").Value = 3309.6
$ws.Range("I83").Value = 3437
$ws.Range("J83").Value = 2800
$ws.Range("K83").Value = 17185
$ws.Range("L83").Value = 14000
$ws.Range("M83").Value = -12193
$ws.Range("N83").Value = -23984
$ws.Range("H112").Value = 15998
$ws.Range("J112").Value = 15998
$ws.Range("L112").Value = 15998
$ws.Range("N112").Value = -18214

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H22").Value = 1945
$ws.Range("I22").Value = 2769.25
$ws.Range("J22").Value = 296.5
$ws.Range("K22").Value = 2769.25
$ws.Range("L22").Value = 296.5
$ws.Range("M22").Value = -2474.25
$ws.Range("N22").Value = -886.5
$ws.Range("H27").Value = 1945
$ws.Range("I27").Value = 2769.25
$ws.Range("J27").Value = 296.5
$ws.Range("K27").Value = 2769.25
$ws.Range("L27").Value = 296.5
$ws.Range("M27").Value = -2662.25
$ws.Range("N27").Value = -510.5
$ws.Range("H136").Value = 35405.64
$ws.Range("I136").Value = 5429.7856
$ws.Range("J136").Value = 111707.82
$ws.Range("K136").Value = 16289.3568
$ws.Range("L136").Value = 335123.46
$ws.Range("M136").Value = -13739.3568
$ws.Range("N136").Value = -340223.46

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H131").Value = 124571.664
$ws.Range("J131").Value = 124571.664
$ws.Range("L131").Value = 124571.664
$ws.Range("N131").Value = -134651.664
